# Clean up the "New HIV infections" sheet: the Number column (B) was stored
# as dirty text like "29000 [22000 - 36000]". Replace each of those with the
# plain numeric value (the first number in the string) so the column holds
# real numbers instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = [string]$cell.Value2

    if ($text -match '^\s*(-?[0-9]+(\.[0-9]+)?)') {
        $numberValue = [double]$matches[1]
        $cell.Value = $numberValue
    }
}
